$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old 6th data row (sheet shrinks from 6 rows to 5 rows)
$ws.Rows.Item(6).Delete()

# Header row
$ws.Range("A1").Value = "Anahtar Kelime"
$ws.Range("B1").Value = "Senaryo"
$ws.Range("C1").Value = "Açıklama"
$ws.Range("D1").Value = "Çözüm"
$ws.Range("E1").Value = "Sorumlu"
$ws.Range("F1").Value = "Görsel"

# Extend the bold/bordered header style from D1 to the new E1:F1 header cells
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Sorumlu"
$ws.Range("F1").Value = "Görsel"

# Row 2
$ws.Range("A2").Value = "dondu"
$ws.Range("B2").Value = "Sistem dondu"
$ws.Range("C2").Value = "Uygulama yanıt vermiyor."
$ws.Range("D2").Value = "Görev yöneticisiyle kapat."
$ws.Range("E2").Value = "BT"
$ws.Range("F2").Value = "dondu_sistem.png"

# Row 3
$ws.Range("A3").Value = "dondu"
$ws.Range("B3").Value = "Bilgisayar dondu"
$ws.Range("C3").Value = "Bilgisayar genel olarak tepki vermiyor."
$ws.Range("D3").Value = "Bilgisayarı yeniden başlat."
$ws.Range("E3").Value = "BT"
$ws.Range("F3").Value = "dondu_pc.png"

# Row 4
$ws.Range("A4").Value = "giriş"
$ws.Range("B4").Value = "Şifre hatası"
$ws.Range("C4").Value = "Kullanıcı adı ya da şifre yanlış girilmiş olabilir."
$ws.Range("D4").Value = "Şifreyi sıfırlayın."
$ws.Range("E4").Value = "Destek"

# Row 5
$ws.Range("A5").Value = "giriş"
$ws.Range("B5").Value = "Hesap kilitli"
$ws.Range("C5").Value = "Çok sayıda yanlış giriş yapılmış olabilir."
$ws.Range("D5").Value = "Destek ekibiyle iletişime geçin."
$ws.Range("E5").Value = "Destek"
